# Cryptos list refresh (prices + 1h volume deltas), plus two row swaps
# (LEO/ImmutableX and Stacks/PEPE traded ranking positions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text in this sheet (some prices use
# thousands-dot grouping like "64.262.56", which is not a valid number).
# Excel auto-converts a plain Range.Value assignment of a numeric-looking
# string into a real number (and would e.g. collapse "0.620" -> 0.62), so
# force text via NumberFormat "@" for the write, then restore the default
# "Normal" style so the cell format matches the rest of the sheet.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '64.262.56'
$ws.Range('E2').Value = '  -4.17%  '

Set-TextValue 'D3' '3.413.24'
$ws.Range('E3').Value = '  -3.36%  '

$ws.Range('E4').Value = '  +0.17%  '

Set-TextValue 'D5' '566.82'
$ws.Range('E5').Value = '  +1.16%  '

Set-TextValue 'D6' '172.76'
$ws.Range('E6').Value = '  -8.42%  '

Set-TextValue 'D7' '0.618'
$ws.Range('E7').Value = '  -0.19%  '

$ws.Range('E8').Value = '  +0.10%  '

Set-TextValue 'D9' '0.620'
$ws.Range('E9').Value = '  -1.89%  '

$ws.Range('E10').Value = '  +2.58%  '

Set-TextValue 'D11' '54.77'
$ws.Range('E11').Value = '  -0.09%  '

Set-TextValue 'D12' '0.0000269'
$ws.Range('E12').Value = '  -0.72%  '

Set-TextValue 'D13' '9.05'
$ws.Range('E13').Value = '  -3.57%  '

Set-TextValue 'D14' '3.971.24'
$ws.Range('E14').Value = '  -2.98%  '

$ws.Range('E15').Value = '  -1.28%  '

Set-TextValue 'D16' '3.414.65'
$ws.Range('E16').Value = '  -3.41%  '

Set-TextValue 'D17' '17.95'
$ws.Range('E17').Value = '  -1.62%  '

Set-TextValue 'D18' '11.78'
$ws.Range('E18').Value = '  -2.42%  '

Set-TextValue 'D19' '64.435.27'
$ws.Range('E19').Value = '  -3.90%  '

Set-TextValue 'D20' '0.985'
$ws.Range('E20').Value = '  -1.25%  '

Set-TextValue 'D21' '404.19'
$ws.Range('E21').Value = '  -5.83%  '

Set-TextValue 'D22' '4.14'
$ws.Range('E22').Value = '  +1.53%  '

Set-TextValue 'D23' '4.38'
$ws.Range('E23').Value = '  +6.23%  '

Set-TextValue 'D24' '82.78'
$ws.Range('E24').Value = '  -3.05%  '

Set-TextValue 'D25' '13.03'
$ws.Range('E25').Value = '  +5.60%  '

Set-TextValue 'D26' '10.75'
$ws.Range('E26').Value = '  -3.04%  '

$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D27' '5.98'
$ws.Range('E27').Value = '  -2.44%  '

$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D28' '2.77'
$ws.Range('E28').Value = '  -4.42%  '

Set-TextValue 'D29' '8.86'
$ws.Range('E29').Value = '  -2.57%  '

Set-TextValue 'D30' '29.60'
$ws.Range('E30').Value = '  -2.71%  '

Set-TextValue 'D31' '6.66'
$ws.Range('E31').Value = '  +1.39%  '

Set-TextValue 'D32' '582.39'
$ws.Range('E32').Value = '  -9.10%  '

Set-TextValue 'D33' '11.43'
$ws.Range('E33').Value = '  -2.55%  '

Set-TextValue 'D34' '0.107'
$ws.Range('E34').Value = '  -3.31%  '

Set-TextValue 'D35' '58.93'
$ws.Range('E35').Value = '  -1.86%  '

$ws.Range('E36').Value = '  +3.88%  '

Set-TextValue 'D37' '0.997'
$ws.Range('E37').Value = '  -0.09%  '

Set-TextValue 'D38' '35.81'
$ws.Range('E38').Value = '  -6.62%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D39' '3.45'
$ws.Range('E39').Value = '  +2.07%  '

$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D40' '0.0₃0753'
$ws.Range('E40').Value = '  -6.89%  '

Set-TextValue 'D41' '0.373'
$ws.Range('E41').Value = '  -4.18%  '

Set-TextValue 'D42' '3.166.59'
$ws.Range('E42').Value = '  +4.05%  '

Set-TextValue 'D43' '1.00'
$ws.Range('E43').Value = '  +0.20%  '

Set-TextValue 'D44' '2.88'
$ws.Range('E44').Value = '  +0.07%  '

Set-TextValue 'D45' '2.48'
$ws.Range('E45').Value = '  -6.21%  '

Set-TextValue 'D46' '3.19'
$ws.Range('E46').Value = '  -4.72%  '

Set-TextValue 'D47' '0.0405'
$ws.Range('E47').Value = '  -3.12%  '

Set-TextValue 'D48' '2.63'
$ws.Range('E48').Value = '  -4.95%  '

Set-TextValue 'D49' '0.129'
$ws.Range('E49').Value = '  -1.79%  '

Set-TextValue 'D50' '8.34'
$ws.Range('E50').Value = '  -3.74%  '

Set-TextValue 'D51' '135.97'
$ws.Range('E51').Value = '  -5.15%  '
